$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.943.79"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.559.75"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.782.27"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.559.88"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "26.950.61"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "0.0₃0704"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").Value = "1.417.76"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.06%  "
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "1.695.96"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0955"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
